# Add the new "Tray-Grade" column header to the bulk tray sheet sample.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header in the next empty column (I1).
$ws.Range("I1").Value = "Tray-Grade"

# Normalise the header row formatting (drop the per-cell style indexes
# that were left over from the old template) so every header cell shares
# the default/Normal style.
$ws.Range("A1:I1").Style = "Normal"

# Leave the selection where the author left it.
$ws.Range("K3").Select() | Out-Null
